$d = $word.ActiveDocument

$replacements = @(
    @("855÷2=", "778÷9="),
    @("739÷9=", "413÷2="),
    @("140÷5=", "754÷3="),
    @("681÷9=", "843÷8="),
    @("896÷6=", "389÷3="),
    @("101÷6=", "684÷7="),
    @("213÷9=", "259÷9="),
    @("107÷7=", "825÷5="),
    @("176÷6=", "905÷8="),
    @("651÷9=", "841÷8="),
    @("660÷3=", "199÷8="),
    @("349÷2=", "319÷4="),
    @("332÷5=", "712÷7="),
    @("618÷8=", "991÷6="),
    @("104÷2=", "372÷3="),
    @("450÷2=", "188÷8="),
    @("872÷4=", "127÷2="),
    @("189÷6=", "543÷9="),
    @("104÷9=", "214÷4="),
    @("713÷3=", "421÷6="),
    @("679÷8=", "395÷8="),
    @("277÷4=", "233÷3="),
    @("139÷3=", "743÷3="),
    @("493÷4=", "138÷2="),
    @("555÷4=", "145÷6=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
